# Revert "Artificial Intelligence and Machine Learning" template content
# back to the "Finance - Core Banking System Modernization" template.
# Also restores the "gap" empty rows that exist in the canonical OOXML
# (rows that are present as bare <row r="N"/> elements with no cell data).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Instructions & User Guide
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Instructions & User Guide")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization Comprehensive Budget - User Guide & Instructions"
$ws.Range("A56").Value = "📋 FINANCE - CORE BANKING SYSTEM MODERNIZATION PROJECT OVERVIEW"
foreach ($r in 10, 20, 28, 37, 45, 54, 55, 60) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Budget Summary
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Budget Summary")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Executive Budget Summary"
foreach ($r in 2, 6) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Resources
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Resources")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Resources Budget"

$ws.Range("A4").Value = "Banking Systems Architect"
$ws.Range("B4").Value = 180
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 14

$ws.Range("A5").Value = "Core Banking Developer"
$ws.Range("B5").Value = 160
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 14

$ws.Range("A6").Value = "Database Administrator"
$ws.Range("B6").Value = 150
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 14

$ws.Range("A7").Value = "Integration Specialist"
$ws.Range("B7").Value = 145
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 14

$ws.Range("A8").Value = "QA/Testing Lead"
$ws.Range("B8").Value = 130
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 14

$ws.Range("A9").Value = "Business Analyst"
$ws.Range("B9").Value = 125
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 14

$ws.Range("A10").Value = "Compliance Officer"
$ws.Range("B10").Value = 140
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 14

foreach ($r in 2, 11) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Logistics
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Logistics")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Logistics Budget"

$ws.Range("B4").Value = 125000
$ws.Range("B5").Value = 180000
$ws.Range("B6").Value = 55000
$ws.Range("B7").Value = 35000
$ws.Range("B8").Value = 25000

foreach ($r in 2, 9) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Technology
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Technology")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Technology Budget"

$ws.Range("A4").Value = "Core Banking Platform License"
$ws.Range("B4").Value = 2500000

$ws.Range("A5").Value = "Cloud Infrastructure (AWS/Azure)"
$ws.Range("B5").Value = 450000

$ws.Range("A6").Value = "Database Management System"
$ws.Range("B6").Value = 280000

$ws.Range("A7").Value = "Security and Encryption Tools"
$ws.Range("B7").Value = 180000

$ws.Range("A8").Value = "API Management Platform"
$ws.Range("B8").Value = 120000

$ws.Range("A9").Value = "Testing and QA Tools"
$ws.Range("B9").Value = 95000

foreach ($r in 2, 10) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Training
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Training")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Training Budget"

$ws.Range("A4").Value = "Core Banking Platform Training"
$ws.Range("B4").Value = 85000
$ws.Range("C4").Value = 45

$ws.Range("A5").Value = "Compliance and Regulatory Training"
$ws.Range("B5").Value = 45000
$ws.Range("C5").Value = 20

$ws.Range("A6").Value = "Technical Skills Development"
$ws.Range("B6").Value = 65000
$ws.Range("C6").Value = 30

$ws.Range("A7").Value = "Change Management Workshops"
$ws.Range("B7").Value = 35000
$ws.Range("C7").Value = 45

$ws.Range("A8").Value = "End-User Training Materials"
$ws.Range("B8").Value = 25000
$ws.Range("C8").Value = 200

foreach ($r in 2, 9) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Contingency
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contingency")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Contingency Budget"

$ws.Range("D6").Value = "Regulatory changes or compliance requirements"
$ws.Range("D8").Value = "Staff turnover or skill gaps in banking domain"
$ws.Range("D9").Value = "Delays or timeline extensions due to testing"

foreach ($r in 2, 5, 11, 13) {
    $ws.Rows.Item($r).OutlineLevel = 0
}

# ---------------------------------------------------------------------
# Sheet: Timeline
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Timeline")
$ws.Range("A1").Value = "Finance - Core Banking System Modernization - Budget Timeline"

foreach ($r in 2) {
    $ws.Rows.Item($r).OutlineLevel = 0
}
